# "Misc changes, progress in excel file"
# Applies the edits captured by the XML diff:
#  - a handful of shared-string (cell text) updates on Sheet1
#  - yellow highlight fill added to the cells that were edited (rows 9, 10, and D:G of row 11)
#  - row-height tweaks on rows 9-11 to match the new (wrapped) content
#  - workbook Date1904 flag explicitly re-asserted (cosmetic attribute-name change upstream)
#  - selection / active-cell state nudged to match the new view

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- workbook-level flag (cosmetic rename of dateCompatibility -> date1904 upstream) ---
$wb.Date1904 = $false

# --- text edits -----------------------------------------------------------
$ws.Range("E9").Value = "Picture of Elephant"
$ws.Range("I9").Value = "Bar chart of coastal borders of Asian countries"

$ws.Range("A10").Value = "Swaziland (Eswatini) only borders one country"
$ws.Range("E10").Value = "Picture of a hockey player"
$ws.Range("F10").Value = "Bar chart of age of residents of Eswatini"
$ws.Range("G10").Value = "Bar chart of age of residents of Canada"

$ws.Range("E11").Value = "Picture of the Nile river"

# --- highlight the edited cells with a yellow fill ------------------------
$yellow = 65535
$ws.Range("D9:N9").Interior.Color = $yellow
$ws.Range("D10:N10").Interior.Color = $yellow
$ws.Range("D11:G11").Interior.Color = $yellow

# --- row heights, re-measured after the text/formatting changes -----------
$ws.Rows.Item(9).RowHeight = 44.95
$ws.Rows.Item(10).RowHeight = 35.6
$ws.Rows.Item(11).RowHeight = 35.6

# --- restore a sensible selection / scroll position ------------------------
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("G7").Select()

Write-Output "done"
